$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text2 = @'
L'analyse des modifications apportées à l'article révèle plusieurs changements notables :
1. **Suppression de certaines obligations spécifiques** : La nouvelle version de l'article a supprimé plusieurs obligations détaillées qui étaient imposées aux prestataires de services, notamment celles relatives à la sécurité et au contrôle interne, à la gestion des conflits d'intérêts, à la résilience et à la sécurité des systèmes informatiques, ainsi qu'aux obligations spécifiques pour les prestataires fournissant le service mentionné au 1° de l'article L. 54-10-2. Ces obligations incluaient la conclusion de conventions avec les clients, l'établissement de politiques de conservation, la restitution des actifs numériques, la ségrégation des détentions et l'usage des actifs numériques avec le consentement des clients.
2. **Réduction des vérifications par l'AMF** : La nouvelle version a également supprimé la vérification par l'Autorité des marchés financiers (AMF) de la sécurité des systèmes d'information des prestataires enregistrés. De plus, la possibilité pour l'AMF de suspendre l'enregistrement d'un prestataire en cas de menace grave et imminente pour la stabilité du marché des actifs numériques a été retirée.
3. **Simplification des procédures** : La nouvelle version semble simplifier les procédures en supprimant certaines étapes et exigences détaillées. Par exemple, la nécessité pour l'AMF de recueillir l'avis conforme de l'Autorité de contrôle prudentiel et de résolution sur le respect des conditions prévues aux 1° à 4° a été maintenue, mais les détails sur les procédures de suspension et de radiation ont été simplifiés.
4. **Clarification des conditions de fourniture de services en France** : La nouvelle version maintient la disposition selon laquelle les conditions dans lesquelles un service est considéré comme fourni en France sont fixées par le règlement général de l'AMF, mais elle a supprimé les détails supplémentaires qui étaient présents dans l'ancienne version.
En résumé, la nouvelle version de l'article semble viser à simplifier et à alléger certaines des obligations et des procédures imposées aux prestataires de services, tout en maintenant les exigences fondamentales en matière d'enregistrement et de conformité. Les modifications apportées pourraient être interprétées comme une volonté de rendre le cadre réglementaire plus flexible et moins contraignant pour les prestataires, tout en conservant les éléments essentiels de supervision et de contrôle par l'AMF.
'@
$ws.Range("V2").Value2 = $text2

$text3 = @'
L'analyse des modifications apportées à cet article révèle plusieurs changements significatifs :
1. **Suppression de la référence au règlement (UE) 2023/1114** : La nouvelle version de l'article ne fait plus mention du règlement (UE) 2023/1114 du Parlement européen et du Conseil du 31 mai 2023 sur les marchés de crypto-actifs. Cela pourrait indiquer une volonté de simplifier le texte en supprimant les références spécifiques à la législation européenne, ou une modification de la portée de l'article pour se concentrer uniquement sur la législation nationale.
2. **Restriction des services concernés** : La nouvelle version de l'article se limite aux services mentionnés aux 1° à 4° de l'article L. 54-10-2, sans inclure les services sur crypto-actifs au sens du règlement européen. Cela pourrait signifier une restriction du champ d'application de l'article, se concentrant uniquement sur les services définis par la législation nationale.
3. **Clarification des interdictions** : La nouvelle version précise que l'interdiction d'utiliser une dénomination, une raison sociale, une publicité ou tout autre procédé laissant croire à un enregistrement ou une autorisation s'applique spécifiquement aux prestataires des services mentionnés aux mêmes 1° et 4°. Cette précision pourrait viser à éviter toute ambiguïté quant aux services concernés par cette interdiction.
En résumé, les modifications apportées à cet article semblent viser à simplifier et clarifier le texte, en se concentrant sur la législation nationale et en précisant les services concernés par les interdictions.
'@
$ws.Range("V3").Value2 = $text3

$text4 = @'
L'analyse des modifications apportées à l'article révèle plusieurs ajustements notables :
1. **Autorité de vérification de la sécurité des systèmes d'information** :
   - **Ancienne version** : L'Autorité des marchés financiers (AMF) pouvait solliciter l'avis de l'autorité mentionnée à l'article L. 2321-1 du code de la défense.
   - **Nouvelle version** : L'AMF peut solliciter l'avis de l'autorité nationale en charge de la sécurité des systèmes d'information. Ce changement clarifie et potentiellement élargit l'autorité consultée pour la sécurité des systèmes d'information.
2. **Seuil de détention de capital ou de droits de vote** :
   - **Ancienne version** : Les personnes physiques détenant plus de 10 % du capital ou des droits de vote du prestataire devaient justifier de leur honorabilité et compétence.
   - **Nouvelle version** : Ce seuil est relevé à plus de 25 %. Cela signifie que les exigences de justification d'honorabilité et de compétence s'appliquent désormais à un nombre potentiellement réduit de détenteurs de capital ou de droits de vote, ce qui pourrait alléger les obligations de certains investisseurs.
3. **Terminologie et précision** :
   - **Ancienne version** : Référence à l'article L. 2321-1 du code de la défense.
   - **Nouvelle version** : Référence plus générale à l'autorité nationale en charge de la sécurité des systèmes d'information. Cette modification peut viser à simplifier et à rendre plus flexible la désignation de l'autorité compétente.
4. **Autres modifications mineures** :
   - Les autres sections de l'article semblent inchangées, à l'exception de quelques ajustements de formulation qui n'affectent pas substantiellement le contenu ou les obligations des prestataires.
En conclusion, les modifications apportées visent principalement à clarifier certaines dispositions, à ajuster les seuils de détention de capital pour les obligations de justification d'honorabilité et de compétence, et à préciser les autorités compétentes en matière de sécurité des systèmes d'information. Ces ajustements peuvent avoir des implications pratiques pour les prestataires de services sur actifs numériques, notamment en termes de conformité et de gouvernance.
'@
$ws.Range("V4").Value2 = $text4

$text5 = @'
L'analyse des modifications apportées à cet article révèle un changement dans la portée des personnes concernées par l'interdiction de tenir des comptes ou livrets d'épargne anonymes. 
L'ancienne version de l'article faisait référence aux personnes mentionnées aux 1° à 7° de l'article L. 561-2. La nouvelle version étend cette référence aux personnes mentionnées aux 1° à 7° bis du même article. 
Cela signifie que la nouvelle version inclut potentiellement des catégories supplémentaires de personnes ou entités qui sont désormais soumises à cette interdiction. Le terme "7° bis" suggère l'ajout d'une nouvelle catégorie ou sous-catégorie dans l'article L. 561-2, élargissant ainsi le champ d'application de la réglementation.
En résumé, la modification vise à inclure davantage de personnes ou entités dans l'interdiction de tenir des comptes ou livrets d'épargne anonymes, renforçant ainsi les mesures de lutte contre le blanchiment d'argent et le financement du terrorisme.
'@
$ws.Range("V5").Value2 = $text5

$text6 = @'
L'analyse des modifications entre l'ancienne et la nouvelle version de l'article révèle plusieurs changements notables :
1. **Suppression de certaines catégories** : 
   - Les prestataires de services sur crypto-actifs ayant leur siège social dans un autre État membre de l'Union européenne ou partie à l'Espace économique européen (ancien 1° quater) ont été supprimés.
   - Les prestataires de services autorisés à fournir des services sur crypto-actifs conformément à l'article 59 du règlement (UE) 2023/1114 (ancien 7° bis b) ont été supprimés.
   - Les conseillers en investissements financiers, dépositaires centraux et sociétés de gestion de placements collectifs autorisés à fournir des services sur crypto-actifs (ancien 7° quater b) ont été supprimés.
   - Les gestionnaires de crédits (ancien 20°) ont été supprimés.
2. **Ajout de nouvelles catégories** :
   - Les émetteurs de jetons ayant obtenu le visa mentionné à l'article L. 552-4 (nouveau 7° ter) ont été ajoutés.
3. **Modifications de terminologie et de références légales** :
   - Les références aux articles du code monétaire et financier ont été mises à jour pour refléter les modifications législatives récentes.
   - Les termes "commissaires de justice" ont été remplacés par "huissiers de justice" et "commissaires-priseurs judiciaires" (nouveau 13°).
4. **Clarifications et précisions** :
   - La nouvelle version apporte des clarifications sur les catégories de personnes assujetties, notamment en ce qui concerne les intermédiaires en financement participatif et les conseillers en investissements participatifs.
   - Les références aux articles spécifiques du code des assurances, du code de la sécurité sociale, du code de la mutualité, et du code de commerce ont été mises à jour pour plus de précision.
5. **Réorganisation et simplification** :
   - La nouvelle version réorganise certaines catégories pour une meilleure lisibilité et simplification, notamment en regroupant les catégories similaires et en supprimant les redondances.
En conclusion, les modifications apportées visent à clarifier, simplifier et mettre à jour les références légales et les catégories de personnes assujetties aux obligations prévues par les dispositions des sections 2 à 7 du chapitre concerné. Ces changements reflètent les évolutions législatives récentes et les besoins de mise en conformité avec les règlements européens.
'@
$ws.Range("V6").Value2 = $text6

$text7 = @'
L'analyse des modifications apportées à l'article révèle plusieurs changements notables :
1. **Extension des entités concernées** :
   - L'ancienne version de l'article mentionnait les personnes énumérées aux 1° à 7° de l'article L. 561-2, ainsi que certaines entités spécifiques comme les compagnies financières holding et les entreprises mères de sociétés de financement. La nouvelle version étend cette liste en incluant les personnes mentionnées aux 1° à 7° bis et au 7° quater de l'article L. 561-2. Cette extension élargit le champ des entités soumises aux obligations de déclaration et de partage d'informations.
2. **Suppression de références spécifiques** :
   - La nouvelle version de l'article a supprimé les références détaillées aux articles spécifiques du code des assurances, du code de la mutualité et du code de la sécurité sociale. Cela simplifie la rédaction de l'article mais pourrait potentiellement réduire la précision des entités concernées.
3. **Ajout d'un nouveau paragraphe (III)** :
   - La nouvelle version introduit un paragraphe III qui permet aux entreprises mères de groupes mentionnées au I de l'article L. 561-33 d'accéder aux informations relatives aux déclarations prévues à l'article L. 561-15 réalisées par les entités du groupe. Cette disposition est soumise à des conditions spécifiques, notamment que les entreprises mères soient établies dans un État membre de l'Union européenne, un État partie à l'accord sur l'Espace économique européen, ou un pays tiers imposant des obligations équivalentes en matière de lutte contre le blanchiment de capitaux et le financement du terrorisme. Cette inclusion renforce la coordination et la supervision au sein des groupes internationaux.
4. **Maintien des conditions de partage d'informations** :
   - Les conditions de partage d'informations (a, b, c, d) restent globalement inchangées dans les deux versions, garantissant que les informations échangées sont nécessaires, utilisées exclusivement à des fins de lutte contre le blanchiment de capitaux et le financement du terrorisme, et que le traitement des informations garantit un niveau de protection suffisant de la vie privée et des libertés et droits fondamentaux des personnes.
5. **Décret en Conseil d'État** :
   - Les deux versions prévoient qu'un décret en Conseil d'État précise les modalités d'application de l'article, assurant ainsi une mise en œuvre réglementaire détaillée et conforme.
En conclusion, les modifications apportées visent principalement à élargir le champ des entités concernées, à simplifier la rédaction de l'article, et à renforcer la coordination au sein des groupes internationaux en matière de lutte contre le blanchiment de capitaux et le financement du terrorisme. Ces changements reflètent une volonté d'améliorer l'efficacité et la cohérence des obligations de déclaration et de partage d'informations au sein des groupes financiers.
'@
$ws.Range("V7").Value2 = $text7

$text8 = @'
L'analyse des modifications apportées à cet article révèle une suppression notable. La nouvelle version de l'article a éliminé la référence spécifique aux prestataires de services sur crypto-actifs au sens du règlement (UE) 2023/1114 du Parlement européen et du Conseil du 31 mai 2023 sur les marchés de crypto-actifs.
Cette suppression pourrait indiquer une volonté de simplifier le texte législatif en retirant une mention spécifique qui pourrait être redondante ou déjà couverte par d'autres dispositions légales. En effet, les prestataires de services sur crypto-actifs pourraient toujours être inclus dans les catégories générales des prestataires de services mentionnés aux 1° à 4° de l'article L. 54-10-2, rendant ainsi la mention spécifique superflue.
En termes de conséquences juridiques, cette modification n'affecte pas les obligations des prestataires de services sur crypto-actifs, qui restent soumis aux mêmes exigences de coopération avec l'Autorité des marchés financiers. Toutefois, la suppression de cette référence spécifique pourrait avoir un impact sur la clarté et la précision de la législation pour les acteurs concernés, nécessitant une interprétation plus large des catégories de prestataires de services mentionnées.
En conclusion, la modification vise probablement à simplifier et à harmoniser le texte législatif sans changer substantiellement les obligations des prestataires de services sur crypto-actifs.
'@
$ws.Range("V8").Value2 = $text8

$text9 = @'
L'analyse des modifications apportées à l'article révèle plusieurs changements significatifs :
1. **Références législatives et réglementaires** :
   - Les articles L. 54-10-3 à L. 54-10-5, initialement régis par la loi n° 2019-486 du 22 mai 2019, sont désormais régis par l'ordonnance n° 2020-1544 du 9 décembre 2020. Cette ordonnance renforce le cadre de la lutte contre le blanchiment de capitaux et le financement du terrorisme applicable aux actifs numériques.
   - L'article L. 572-24, initialement sous la loi n° 2019-486, est également régi par l'ordonnance n° 2020-1544.
2. **Adaptations spécifiques pour la Nouvelle-Calédonie** :
   - Pour l'application de l'article L. 54-10-3, certaines mentions spécifiques à l'Union européenne sont supprimées, ce qui indique une adaptation locale des dispositions européennes.
   - Les références aux règlements européens dans les articles L. 54-10-3 et L. 54-10-5 sont remplacées par des références locales, spécifiquement à l'article L. 713-16.
   - Les références au code de commerce sont remplacées par des dispositions équivalentes applicables localement, ce qui montre une volonté d'adapter les textes aux réalités juridiques locales.
3. **Conversion des montants financiers** :
   - Les montants en euros mentionnés dans les articles L. 572-23 et L. 572-26 sont convertis en francs CFP, la monnaie utilisée en Nouvelle-Calédonie. Par exemple, 15 000 euros deviennent 1 790 000 francs CFP, et 30 000 euros deviennent 3 580 000 francs CFP. Cette conversion assure la pertinence et l'applicabilité des sanctions financières dans le contexte local.
4. **Suppression et remplacement de termes** :
   - Les mots “la Banque de France” sont remplacés par “l'Institut d'émission d'outre-mer” pour l'application du dernier alinéa de l'article L. 54-10-5, ce qui est maintenu dans la nouvelle version.
En résumé, les modifications apportées visent principalement à adapter les dispositions législatives et réglementaires aux spécificités locales de la Nouvelle-Calédonie, en tenant compte des nouvelles régulations introduites par l'ordonnance de 2020 et en assurant la conversion des montants financiers dans la monnaie locale. Ces ajustements permettent une application plus cohérente et efficace des textes de loi dans ce territoire d'outre-mer.
'@
$ws.Range("V9").Value2 = $text9

$text10 = @'
L'analyse des modifications entre l'ancienne et la nouvelle version de l'article révèle plusieurs changements notables :
1. **Références législatives et ordonnances** :
   - La nouvelle version a supprimé les références spécifiques à certaines lois et ordonnances, notamment celles de la loi n° 2021-1308 du 8 octobre 2021 et de l'ordonnance n° 2021-958 du 19 juillet 2021. Cela simplifie la lecture et l'application de l'article en se concentrant sur les ordonnances et lois les plus pertinentes.
2. **Applicabilité des articles** :
   - Les articles L. 561-2, L. 561-7, L. 561-14 et L. 561-20 sont désormais applicables dans leur rédaction résultant de l'ordonnance n° 2020-1544 du 9 décembre 2020, ce qui unifie et clarifie les références législatives.
3. **Suppression de dispositions spécifiques** :
   - La nouvelle version a supprimé le point 7° du III, ce qui pourrait indiquer une simplification ou une suppression de dispositions jugées non nécessaires ou redondantes.
4. **Adaptations locales** :
   - Les adaptations locales restent largement inchangées, mais la nouvelle version a supprimé certaines références spécifiques aux règlements européens, ce qui pourrait indiquer une volonté de simplifier l'application locale en Nouvelle-Calédonie.
5. **Clarifications et simplifications** :
   - La nouvelle version a supprimé plusieurs phrases et références spécifiques, ce qui semble viser à simplifier et clarifier le texte pour une meilleure compréhension et application.
En résumé, les modifications apportées visent principalement à simplifier et clarifier l'article, en supprimant des références spécifiques et en unifiant les textes applicables. Cela devrait faciliter l'application des dispositions en Nouvelle-Calédonie, tout en maintenant les adaptations nécessaires pour tenir compte des spécificités locales.
'@
$ws.Range("V10").Value2 = $text10

$text11 = @'
L'analyse des modifications entre l'ancienne et la nouvelle version de l'article révèle plusieurs changements notables :
1. **Références législatives et ordonnances** :
   - La nouvelle version a supprimé les références spécifiques à certaines lois et ordonnances, notamment celles de la loi n° 2021-1308 du 8 octobre 2021 et de l'ordonnance n° 2021-958 du 19 juillet 2021. Cela simplifie la rédaction et peut indiquer une volonté de se concentrer sur des textes plus récents ou plus pertinents.
2. **Applicabilité des articles** :
   - Les articles L. 561-2, L. 561-7, L. 561-14 et L. 561-20 sont désormais applicables dans leur rédaction résultant de l'ordonnance n° 2020-1544 du 9 décembre 2020, ce qui montre une mise à jour pour inclure des dispositions plus récentes concernant la lutte contre le blanchiment de capitaux et le financement du terrorisme.
3. **Suppression de certaines dispositions** :
   - Le point 7° du III a été supprimé dans la nouvelle version. Cela pourrait indiquer que les dispositions précédemment mentionnées ne sont plus jugées nécessaires ou pertinentes pour l'application en Polynésie française.
4. **Clarification et simplification** :
   - La nouvelle version semble avoir simplifié certaines formulations et supprimé des redondances, ce qui peut faciliter la compréhension et l'application des dispositions légales.
5. **Adaptations locales** :
   - Les adaptations locales restent largement inchangées, mais la suppression de certaines références européennes et la clarification des termes montrent une volonté d'adapter plus précisément les textes aux réalités locales de la Polynésie française.
En résumé, les modifications apportées visent principalement à actualiser les références législatives, simplifier la rédaction, et supprimer des dispositions jugées non nécessaires, tout en maintenant une adaptation locale cohérente. Ces changements reflètent une volonté d'aligner les textes sur les évolutions législatives récentes et de faciliter leur application en Polynésie française.
'@
$ws.Range("V11").Value2 = $text11

